$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 3100.3635
$ws.Range("I40").Value = 8102
$ws.Range("J40").Value = 1988.8889
$ws.Range("K40").Value = 8102
$ws.Range("L40").Value = 1988.8889
$ws.Range("M40").Value = -7927
$ws.Range("N40").Value = -2338.8889
$ws.Range("H51").Value = 4066.6667
$ws.Range("J51").Value = 4066.6667
$ws.Range("L51").Value = 4066.6667
$ws.Range("N51").Value = -5034.6667
$ws.Range("H106").Value = 2004201
$ws.Range("I106").Value = 1255251.2
$ws.Range("J106").Value = 5000000
$ws.Range("K106").Value = 1255251.2
$ws.Range("L106").Value = 5000000
$ws.Range("M106").Value = -1254620.2
$ws.Range("N106").Value = -5001262
$ws.Range("H113").Value = 106036.86
$ws.Range("J113").Value = 4982.857
$ws.Range("L113").Value = 4982.857
$ws.Range("N113").Value = -11490.857
$ws.Range("H137").Value = 5717346.5
$ws.Range("I137").Value = 1684.1
$ws.Range("J137").Value = 13338229
$ws.Range("K137").Value = 5052.299999999999
$ws.Range("L137").Value = 40014687
$ws.Range("M137").Value = -2502.299999999999
$ws.Range("N137").Value = -40019787

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8274.528
$ws.Range("I32").Value = 9187.647999999999
$ws.Range("J32").Value = 5192.75
$ws.Range("K32").Value = 9187.647999999999
$ws.Range("L32").Value = 5192.75
$ws.Range("M32").Value = -8900.647999999999
$ws.Range("N32").Value = -5766.75
$ws.Range("H61").Value = 11907334
$ws.Range("I61").Value = 16668947
$ws.Range("J61").Value = 3303.5
$ws.Range("K61").Value = 16668947
$ws.Range("L61").Value = 3303.5
$ws.Range("M61").Value = -16668735
$ws.Range("N61").Value = -3727.5
$ws.Range("H97").Value = 7481.4707
$ws.Range("I97").Value = 12664.444
$ws.Range("J97").Value = 1650.625
$ws.Range("K97").Value = 12664.444
$ws.Range("L97").Value = 1650.625
$ws.Range("M97").Value = -12168.444
$ws.Range("N97").Value = -2642.625
$ws.Range("H132").Value = 7814933.5
$ws.Range("I132").Value = 13159781
$ws.Range("J132").Value = 3233.4614
$ws.Range("K132").Value = 39479343
$ws.Range("L132").Value = 9700.3842
$ws.Range("M132").Value = -39476813
$ws.Range("N132").Value = -14760.3842
$ws.Range("H136").Value = 11907334
$ws.Range("I136").Value = 16668947
$ws.Range("J136").Value = 3303.5
$ws.Range("K136").Value = 50006841
$ws.Range("L136").Value = 9910.5
$ws.Range("M136").Value = -50004291
$ws.Range("N136").Value = -15010.5

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 5408.0444
$ws.Range("I134").Value = 4690.5713
$ws.Range("J134").Value = 6589.7646
$ws.Range("K134").Value = 14071.7139
$ws.Range("L134").Value = 19769.2938
$ws.Range("M134").Value = -11536.7139
$ws.Range("N134").Value = -24839.2938

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 2500
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 2500
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 2500
$ws.Range("M62").ClearContents()
$ws.Range("N62").Value = -3748
$ws.Range("H65").Value = 2500
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 2500
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 12500
$ws.Range("M65").ClearContents()
$ws.Range("N65").Value = -18740
$ws.Range("H93").Value = 9313.333000000001
$ws.Range("I93").Value = 9313.333000000001
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 9313.333000000001
$ws.Range("L93").Value = 0
$ws.Range("M93").Value = -7441.333000000001
$ws.Range("N93").ClearContents()
$ws.Range("H96").Value = 24648.334
$ws.Range("J96").Value = 24648.334
$ws.Range("L96").Value = 24648.334
$ws.Range("N96").Value = -30140.334
$ws.Range("H114").Value = 24934.8
$ws.Range("J114").Value = 24934.8
$ws.Range("L114").Value = 24934.8
$ws.Range("N114").Value = -33612.8
$ws.Range("H134").Value = 542959.75
$ws.Range("I134").Value = 1865.3478
$ws.Range("J134").Value = 1135586.9
$ws.Range("K134").Value = 5596.0434
$ws.Range("L134").Value = 3406760.7
$ws.Range("M134").Value = -3061.0434
$ws.Range("N134").Value = -3411830.7

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 5579.8125
$ws.Range("I3").Value = 4338.3335
$ws.Range("J3").Value = 7176
$ws.Range("K3").Value = 13015.0005
$ws.Range("L3").Value = 21528
$ws.Range("M3").Value = -12903.0005
$ws.Range("N3").Value = -21752
$ws.Range("H131").Value = 826.7
$ws.Range("J131").Value = 897.7528
$ws.Range("L131").Value = 2693.2584
$ws.Range("N131").Value = -12773.2584
$ws.Range("H132").Value = 2346.074
$ws.Range("J132").Value = 3599.5715
$ws.Range("L132").Value = 32396.1435
$ws.Range("N132").Value = -37456.1435
$ws.Range("H137").Value = 5055387.5
$ws.Range("I137").Value = 12823641
$ws.Range("J137").Value = 6022.95
$ws.Range("K137").Value = 38470923
$ws.Range("L137").Value = 18068.85
$ws.Range("M137").Value = -38465823
$ws.Range("N137").Value = -28268.85
$ws.Range("H139").Value = 1942
$ws.Range("I139").Value = 1222.3334
$ws.Range("K139").Value = 3667.0002
$ws.Range("M139").Value = 1472.9998
$ws.Range("H140").Value = 3476.8708
$ws.Range("I140").Value = 1483.3684
$ws.Range("K140").Value = 4450.1052
$ws.Range("M140").Value = 729.8948
$ws.Range("H141").Value = 7374.3
$ws.Range("I141").Value = 7177.5
$ws.Range("J141").Value = 7505.5
$ws.Range("K141").Value = 21532.5
$ws.Range("L141").Value = 22516.5
$ws.Range("M141").Value = -16352.5
$ws.Range("N141").Value = -32876.5

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 2000
$ws.Range("I113").Value = 2000
$ws.Range("K113").Value = 2000
$ws.Range("M113").Value = 170
$ws.Range("H126").Value = 4777.0586
$ws.Range("I126").Value = 2404
$ws.Range("J126").Value = 5285.5713
$ws.Range("K126").Value = 7212
$ws.Range("L126").Value = 15856.7139
$ws.Range("M126").Value = -4742
$ws.Range("N126").Value = -20796.7139
$ws.Range("H132").Value = 6161.1055
$ws.Range("I132").Value = 5067
$ws.Range("K132").Value = 15201
$ws.Range("M132").Value = -12671
$ws.Range("H133").Value = 54874.5
$ws.Range("J133").Value = 54874.5
$ws.Range("L133").Value = 54874.5
$ws.Range("N133").Value = -64994.5

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1511.375
$ws.Range("I22").Value = 1020
$ws.Range("J22").Value = 1734.7273
$ws.Range("K22").Value = 1020
$ws.Range("L22").Value = 1734.7273
$ws.Range("M22").Value = -725
$ws.Range("N22").Value = -2324.7273
$ws.Range("H27").Value = 1511.375
$ws.Range("I27").Value = 1020
$ws.Range("J27").Value = 1734.7273
$ws.Range("K27").Value = 1020
$ws.Range("L27").Value = 1734.7273
$ws.Range("M27").Value = -913
$ws.Range("N27").Value = -1948.7273
$ws.Range("H46").Value = 1073.375
$ws.Range("I46").Value = 834.5
$ws.Range("K46").Value = 834.5
$ws.Range("M46").Value = -646.5
$ws.Range("H55").Value = 269.12
$ws.Range("J55").Value = 385.85715
$ws.Range("L55").Value = 385.85715
$ws.Range("N55").Value = -731.85715
$ws.Range("H93").Value = 1299.2273
$ws.Range("I93").Value = 1115.6154
$ws.Range("K93").Value = 1115.6154
$ws.Range("M93").Value = 132.3846000000001

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H74").Value = 10275
$ws.Range("J74").Value = 10275
$ws.Range("L74").Value = 10275
$ws.Range("N74").Value = -12147
$ws.Range("H77").Value = 10275
$ws.Range("J77").Value = 10275
$ws.Range("L77").Value = 30825
$ws.Range("N77").Value = -40185
$ws.Range("H96").Value = 4089.1
$ws.Range("I96").Value = 2998.6667
$ws.Range("J96").Value = 4556.4287
$ws.Range("K96").Value = 2998.6667
$ws.Range("L96").Value = 4556.4287
$ws.Range("M96").Value = -1625.6667
$ws.Range("N96").Value = -7302.4287
